$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = '''64.023.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '''2.757.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''575.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.58%  '
$ws.Range("D6").Value = '''159.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("D8").Value = '''0.601'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.74%  '
$ws.Range("E9").Value = '  -3.81%  '
$ws.Range("E10").Value = '  +3.96%  '
$ws.Range("D11").Value = '''5.81'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -15.42%  '
$ws.Range("D12").Value = '''0.387'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.98%  '
$ws.Range("D13").Value = '''3.247.16'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").Value = '''26.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.28%  '
$ws.Range("D15").Value = '''63.658.55'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("E16").Value = '  -5.45%  '
$ws.Range("D17").Value = '''2.762.74'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '''12.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("E19").Value = '  -3.19%  '
$ws.Range("D20").Value = '''356.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.26%  '
$ws.Range("E21").Value = '  -5.68%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").Value = '''0.528'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.62%  '
$ws.Range("D24").Value = '''65.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.88%  '
$ws.Range("E25").Value = '  -2.20%  '
$ws.Range("D26").Value = '''8.59'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.36%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '''0.0₃0904'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.50%  '
$ws.Range("D29").Value = '''7.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.46%  '
$ws.Range("E30").Value = '  -4.91%  '
$ws.Range("D31").Value = '''1.26'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").Value = '''169.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.09%  '
$ws.Range("D33").Value = '''4.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.42%  '
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("E35").Value = '  -0.34%  '
$ws.Range("D36").Value = '''0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("E37").Value = '  -2.73%  '
$ws.Range("E38").Value = '  -2.20%  '
$ws.Range("D39").Value = '''349.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.80%  '
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("D42").Value = '''39.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("E45").Value = '  -3.34%  '
$ws.Range("E48").Value = '  -2.82%  '
$ws.Range("D49").Value = '''0.101'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.62%  '
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("D51").Value = '''11.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.04%  '

# --- Row swaps (Coin / Link / Price / Volume) for rows 43/44 and 46/47 ---
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '''21.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.17%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''21.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.87%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''136.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.12%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '''0.633'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.35%  '
